# Hjemme passive updated meanEMG legmaxROM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - subject/trial identifiers for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - meanEMG values updated; C2 and E2 no longer have data
$ws.Range("B2").Value = 28.55068550296312
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 18.210526618971528
$ws.Range("E2").ClearContents()

# Row 3 - legmaxROM values updated
$ws.Range("B3").Value = 23.738576534686445
$ws.Range("C3").Value = 23.27312537304482
$ws.Range("D3").Value = 17.814424024323046
$ws.Range("E3").Value = 35.843641919456275

# Update the active selection to reflect the updated data range
$ws.Range("B1:E3").Select()
